$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# --- Simple value updates (rows keep their position) ---
$ws.Range("B3").Value  = "0.1.7"                                    # Version
$ws.Range("B6").Value  = "draft"                                    # Status
$ws.Range("B8").Value  = "2024-11-22T12:33:30-06:00"                # Date
$ws.Range("B10").Value = "The Medical College of Wisconsin, Inc. and the National Marrow Donor Program (http://www.cibmtr.org)"  # Contact (org)
$ws.Range("B11").Value = "Bob Milius (bmilius@nmdp.org)"            # Contact (person)

# --- Insert a new "Jurisdiction" row after the Contact rows (row 12) ---
# Copy the formatting the new trailing row needs (row 16 is brand new) before
# shifting values, so every row keeps the same style index (s="2") as before.
$ws.Range("A15:B15").Copy()
$ws.Range("A16:B16").PasteSpecial(-4122)

# Shift rows 12-15 down to 13-16 (bottom-up so we don't clobber source values)
$ws.Range("A16").Value = $ws.Range("A15").Value()
$ws.Range("B16").Value = $ws.Range("B15").Value()

$ws.Range("A15").Value = $ws.Range("A14").Value()
$ws.Range("B15").Value = $ws.Range("B14").Value()

$ws.Range("A14").Value = $ws.Range("A13").Value()
$ws.Range("B14").Value = $ws.Range("B13").Value()

$ws.Range("A13").Value = $ws.Range("A12").Value()
$ws.Range("B13").Value = $ws.Range("B12").Value()

# New row 12 content
$ws.Range("A12").Value = "Jurisdiction"
$ws.Range("B12").Value = ""
